# Scheduled-runner price refresh for the Leve profit tracking sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Updates currentAveragePrice*
# and derived Leve profit columns (H:N) from the latest market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 4916.5
$ws.Range("J58").Value = 4916.5
$ws.Range("L58").Value = 14749.5
$ws.Range("N58").Value = -15049.5
$ws.Range("H70").Value = 87487.91
$ws.Range("I70").Value = 225254.75
$ws.Range("J70").Value = 2708.3076
$ws.Range("K70").Value = 675764.25
$ws.Range("L70").Value = 8124.9228
$ws.Range("M70").Value = -675494.25
$ws.Range("N70").Value = -8664.9228
$ws.Range("H73").Value = 87487.91
$ws.Range("I73").Value = 225254.75
$ws.Range("J73").Value = 2708.3076
$ws.Range("K73").Value = 675764.25
$ws.Range("L73").Value = 8124.9228
$ws.Range("M73").Value = -674828.25
$ws.Range("N73").Value = -9996.9228
$ws.Range("H112").Value = 1811.1875
$ws.Range("J112").Value = 2057.182
$ws.Range("L112").Value = 6171.545999999999
$ws.Range("N112").Value = -8387.545999999998
$ws.Range("H113").Value = 6207.353
$ws.Range("I113").Value = 4494.1816
$ws.Range("J113").Value = 9348.166999999999
$ws.Range("K113").Value = 4494.1816
$ws.Range("L113").Value = 9348.166999999999
$ws.Range("M113").Value = -1240.1816
$ws.Range("N113").Value = -15856.167
$ws.Range("H132").Value = 304329.94
$ws.Range("I132").Value = 1410.2693
$ws.Range("K132").Value = 4230.8079
$ws.Range("M132").Value = -1700.8079
$ws.Range("H137").Value = 3058.25
$ws.Range("I137").Value = 2657.2942
$ws.Range("K137").Value = 7971.882599999999
$ws.Range("M137").Value = -5421.882599999999
$ws.Range("H138").Value = 2847.2598
$ws.Range("J138").Value = 2462.5688
$ws.Range("L138").Value = 7387.7064
$ws.Range("N138").Value = -17667.7064

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1313.25
$ws.Range("I45").Value = 1351.25
$ws.Range("K45").Value = 1351.25
$ws.Range("M45").Value = -974.25
$ws.Range("H61").Value = 2427.606
$ws.Range("I61").Value = 2358.3845
$ws.Range("K61").Value = 2358.3845
$ws.Range("M61").Value = -2146.3845
$ws.Range("H88").Value = 3179.2
$ws.Range("I88").Value = 2665.3333
$ws.Range("J88").Value = 3950
$ws.Range("K88").Value = 2665.3333
$ws.Range("L88").Value = 3950
$ws.Range("M88").Value = -2259.3333
$ws.Range("N88").Value = -4762
$ws.Range("H91").Value = 3179.2
$ws.Range("I91").Value = 2665.3333
$ws.Range("J91").Value = 3950
$ws.Range("K91").Value = 2665.3333
$ws.Range("L91").Value = 3950
$ws.Range("M91").Value = -1261.3333
$ws.Range("N91").Value = -6758
$ws.Range("H110").Value = 1572.7894
$ws.Range("J110").Value = 1746.8334
$ws.Range("L110").Value = 1746.8334
$ws.Range("N110").Value = -5836.8334
$ws.Range("H122").Value = 2726.6667
$ws.Range("I122").Value = 2790.8
$ws.Range("J122").Value = 2566.3333
$ws.Range("K122").Value = 8372.400000000001
$ws.Range("L122").Value = 7698.999899999999
$ws.Range("M122").Value = -5922.400000000001
$ws.Range("N122").Value = -12598.9999
$ws.Range("H136").Value = 2427.606
$ws.Range("I136").Value = 2358.3845
$ws.Range("K136").Value = 7075.1535
$ws.Range("M136").Value = -4525.1535

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 828.38464
$ws.Range("I107").Value = 730.75
$ws.Range("K107").Value = 730.75
$ws.Range("M107").Value = 1189.25
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1646.2195
$ws.Range("I31").Value = 1218.5714
$ws.Range("K31").Value = 1218.5714
$ws.Range("M31").Value = -923.5714
$ws.Range("H34").Value = 1646.2195
$ws.Range("I34").Value = 1218.5714
$ws.Range("K34").Value = 1218.5714
$ws.Range("M34").Value = -1016.5714
$ws.Range("H58").Value = 3399.3
$ws.Range("I58").Value = 3332.5557
$ws.Range("K58").Value = 3332.5557
$ws.Range("M58").Value = -3129.5557
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H86").Value = 5004.3335
$ws.Range("I86").Value = 5007
$ws.Range("J86").Value = 5003
$ws.Range("K86").Value = 5007
$ws.Range("L86").Value = 5003
$ws.Range("M86").Value = -3884
$ws.Range("N86").Value = -7249
$ws.Range("H89").Value = 5004.3335
$ws.Range("I89").Value = 5007
$ws.Range("J89").Value = 5003
$ws.Range("K89").Value = 25035
$ws.Range("L89").Value = 25015
$ws.Range("M89").Value = -19419
$ws.Range("N89").Value = -36247
$ws.Range("H132").Value = 2512.3333
$ws.Range("I132").Value = 2011.5
$ws.Range("J132").Value = 3514
$ws.Range("K132").Value = 6034.5
$ws.Range("L132").Value = 10542
$ws.Range("M132").Value = -3504.5
$ws.Range("N132").Value = -15602
$ws.Range("H136").Value = 3399.3
$ws.Range("I136").Value = 3332.5557
$ws.Range("K136").Value = 9997.667099999999
$ws.Range("M136").Value = -7447.667099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1160.5834
$ws.Range("J5").Value = 993
$ws.Range("L5").Value = 2979
$ws.Range("N5").Value = -3203
$ws.Range("H113").Value = 894.44446
$ws.Range("I113").Value = 702
$ws.Range("K113").Value = 2106
$ws.Range("M113").Value = 64
$ws.Range("H132").Value = 1513
$ws.Range("I132").Value = 1884.6666
$ws.Range("K132").Value = 16961.9994
$ws.Range("M132").Value = -14431.9994
$ws.Range("H135").Value = 1160.5834
$ws.Range("J135").Value = 993
$ws.Range("L135").Value = 8937
$ws.Range("N135").Value = -14007

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9176.154
$ws.Range("I70").Value = 7782.4
$ws.Range("J70").Value = 10047.25
$ws.Range("K70").Value = 7782.4
$ws.Range("L70").Value = 10047.25
$ws.Range("M70").Value = -7512.4
$ws.Range("N70").Value = -10587.25
$ws.Range("H73").Value = 9176.154
$ws.Range("I73").Value = 7782.4
$ws.Range("J73").Value = 10047.25
$ws.Range("K73").Value = 7782.4
$ws.Range("L73").Value = 10047.25
$ws.Range("M73").Value = -6846.4
$ws.Range("N73").Value = -11919.25
$ws.Range("H102").Value = 3002.625
$ws.Range("I102").Value = 2220.3333
$ws.Range("K102").Value = 2220.3333
$ws.Range("M102").Value = -598.3332999999998
$ws.Range("H132").Value = 2239.2307
$ws.Range("I132").Value = 2114.375
$ws.Range("J132").Value = 2439
$ws.Range("K132").Value = 6343.125
$ws.Range("L132").Value = 7317
$ws.Range("M132").Value = -3813.125
$ws.Range("N132").Value = -12377

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4339.607
$ws.Range("I122").Value = 4748
$ws.Range("J122").Value = 4075.353
$ws.Range("K122").Value = 14244
$ws.Range("L122").Value = 12226.059
$ws.Range("M122").Value = -11794
$ws.Range("N122").Value = -17126.059
$ws.Range("H136").Value = 4098.7
$ws.Range("I136").Value = 4853
$ws.Range("K136").Value = 14559
$ws.Range("M136").Value = -12009

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2219.75
$ws.Range("I81").Value = 2219.75
$ws.Range("K81").Value = 4439.5
$ws.Range("M81").Value = -3378.5
$ws.Range("H84").Value = 2219.75
$ws.Range("I84").Value = 2219.75
$ws.Range("K84").Value = 22197.5
$ws.Range("M84").Value = -16893.5
$ws.Range("H132").Value = 1305.7333
$ws.Range("J132").Value = 794.6667
$ws.Range("L132").Value = 2384.0001
$ws.Range("N132").Value = -7444.0001
